$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "Polygon"
$ws.Range("C9").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D9").Value = "'1.119"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.82%  "

$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "'40.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -14.12%  "

$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'1.660"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -16.27%  "

$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "'0.08145"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.10%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06160"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.54%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.214"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.17%  "

$ws.Range("D2").Value = "21.995.68"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "1.552.98"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D6").Value = "'286.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("D7").Value = "'0.3776"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.17%  "
$ws.Range("D8").Value = "'0.3224"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.70%  "
$ws.Range("D11").Value = "'0.07282"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.98%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "'19.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.04%  "
$ws.Range("D14").Value = "'5.705"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.89%  "
$ws.Range("D15").Value = "'6.767"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").Value = "1.552.13"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "'0.00001079"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.16%  "
$ws.Range("D18").Value = "'0.06631"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").Value = "'84.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.84%  "
$ws.Range("D20").Value = "'6.420"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "'15.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'11.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.63%  "
$ws.Range("D24").Value = "22.010.07"
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("D25").Value = "'2.262"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.60%  "
$ws.Range("D26").Value = "'2.510"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.73%  "
$ws.Range("D27").Value = "'148.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").Value = "'18.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.20%  "
$ws.Range("D29").Value = "'4.853"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("D30").Value = "1.725.62"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("D31").Value = "'120.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.70%  "
$ws.Range("D32").Value = "'1.111"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").Value = "'5.896"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.08%  "
$ws.Range("D36").Value = "'9.253"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.22%  "
$ws.Range("D39").Value = "'0.02271"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.89%  "
$ws.Range("D40").Value = "'0.2110"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.62%  "
$ws.Range("D41").Value = "'1.212"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.29%  "
$ws.Range("D42").Value = "'10.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.82%  "
$ws.Range("D44").Value = "'0.5922"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.02%  "
$ws.Range("D45").Value = "'13.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.57%  "
$ws.Range("D46").Value = "'3.714"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").Value = "'0.5710"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.35%  "
$ws.Range("D48").Value = "'119.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.31%  "
$ws.Range("D49").Value = "'1.922"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.78%  "
$ws.Range("D50").Value = "'1.153"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.46%  "
$ws.Range("D51").Value = "'0.06887"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.29%  "
